# BOOK-16: implementation of job logic and fix of old classes what were
# involved into book saving process.
#
# Appends three new 8-volume book collections (Жюль Верн / Майн Рид /
# Антон Павлович Чехов) to "Лист1", right after the existing data that
# ends on row 179.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

$authors = @(
    @{ Name = "Жюль Верн";             Title = "Собрание сочинений в 8-ти томах";  Year = 1985 },
    @{ Name = "Майн Рид";              Title = "Собрание сочинений в 12-ти томах"; Year = 1992 },
    @{ Name = "Антон Павлович Чехов";  Title = "Собрание сочинений в 8-ти томах";  Year = 1970 }
)

$row = 180
foreach ($author in $authors) {
    for ($volume = 1; $volume -le 8; $volume++) {
        $ws.Cells.Item($row, 1).Value = $author.Name
        $ws.Cells.Item($row, 2).Value = $author.Title
        $ws.Cells.Item($row, 3).Value = $volume
        $ws.Cells.Item($row, 4).Value = 8
        $ws.Cells.Item($row, 5).Value = $author.Year
        $row++
    }
}

# Move the visible selection to where the user left off entering data.
$ws.Range("B204").Select()
